$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($row, $values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# Row 10 - Rusia (values update only)
Set-Row 10 @(145268,10581,18095,125817,2300,76,1356)

# Row 36 - Polonia (values update only)
Set-Row 36 @(13937,244,4095,9159,160,5,683)

# Row 41 - Dinamarca (values update only)
Set-Row 41 @(9670,147,6987,2199,62,0,484)

# Rows 42-44 - country labels rotate (Serbia->Filipinas row42, Banglades->Serbia row43, Filipinas->Banglades row44)
# along with refreshed statistics for each.
$ws.Cells.Item(42, 1).Value = "Filipinas"
Set-Row 42 @(9485,262,1315,7547,31,16,623)

$ws.Cells.Item(43, 1).Value = "Serbia"
Set-Row 43 @(9464,0,1551,7720,54,0,193)

$ws.Cells.Item(44, 1).Value = "Banglades"
Set-Row 44 @(9455,0,1063,8215,1,0,177)

# Row 46 - Noruega (values update only)
Set-Row 46 @(7847,0,32,7602,37,2,213)

# Row 60 - Kazajistan (values update only)
Set-Row 60 @(3988,68,1084,2877,40,0,27)

# Rows 66-67 - country labels swap (Grecia->Oman row66, Oman->Grecia row67) with refreshed stats
$ws.Cells.Item(66, 1).Value = "Oman"
Set-Row 66 @(2637,69,816,1809,17,0,12)

$ws.Cells.Item(67, 1).Value = "Grecia"
Set-Row 67 @(2626,0,1374,1108,37,0,144)

# Row 72 - Uzbekistan (values update only)
Set-Row 72 @(2160,11,1327,823,8,0,10)

# Row 78 - Estonia (values update only)
Set-Row 78 @(1703,3,259,1389,6,0,55)

# Row 87 - Eslovaquia (values update only)
Set-Row 87 @(1413,5,643,745,7,1,25)

# Row 144 - Brunei (values update only)
Set-Row 144 @(138,0,130,7,2,0,1)
